$d = $word.ActiveDocument

# --- 1. Fix the two typos in the summary paragraph -----------------------
# "This section provide a summary ... mesurement " ->
# "This section provides a summary ... measurement "
# (Word's Find/Replace merges every run in the paragraph that shares the
#  same character formatting into a single run, so we restore the original
#  run boundaries afterwards using zero-width bookmarks, which split a run
#  without merging anything else.)
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute( `
    "This section provide a summary of the most severe security vulnerability identified in the structural quality analysis and mesurement ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "This section provides a summary of the most severe security vulnerability identified in the structural quality analysis and measurement ", `
    2)

# --- 2. Re-split the (now merged) run back into its original pieces, plus
#        the two new pieces introduced by the typo fixes.
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$findRng.Find.Execute("This section provides", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $findRng.Start

# Offsets (relative to $base) of every run boundary in the target text:
#   "This section provide"                                                  -> 20
#   "s"                                                                     -> 21
#   " a summary of the most severe s"                                       -> 52
#   "e"                                                                     -> 53
#   "curity vulnerability identified in the structural quality analysis and me" -> 126
#   "a"                                                                     -> 127
#   [ _GoBack bookmark ]                                                    -> 136
#   "surement "                                                             -> 136
#   "by"                                                                    -> 138
#   " CAST AIP"                                                             -> 147
#   " against the "                                                         -> 160
$splitOffsets = @(20, 21, 52, 53, 126, 127, 136, 138, 147, 160)
foreach ($off in $splitOffsets) {
    $pos = $base + $off
    $splitPoint = $d.Range($pos, $pos)
    $d.Bookmarks.Add("ZZZTempSplit", $splitPoint)
}

# --- 3. Re-create the "_GoBack" bookmark at the point where the last edit
#        landed (right after "...analysis and mea", before "surement "),
#        matching the target revision. The old "_GoBack" that used to sit
#        in the Applicability table further down is removed in step 4.
$goBackPos = $base + 127
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
